# Updates the crypto price/volume table per the scraped-data refresh.
# Some "Price" values are numeric-looking strings (e.g. "0.999", "8.93");
# assigning them directly would make Excel auto-convert the cell to a
# number. To preserve the original text-cell representation, we force a
# text number format before assigning, then restore the cell style so no
# stray formatting attribute is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "69.327.52"
$ws.Range("E2").Value = "  -0.20%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.673.05"
$ws.Range("E3").Value = "  -0.47%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "646.04"
$ws.Range("E5").Value = "  -5.23%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "158.18"
$ws.Range("E6").Value = "  -1.25%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - XRP
Set-TextValue $ws.Range("D8") "0.498"
$ws.Range("E8").Value = "  +0.27%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -1.37%  "

# Row 10 - Toncoin
Set-TextValue $ws.Range("D10") "7.09"
$ws.Range("E10").Value = "  -0.52%  "

# Row 11 - Cardano
Set-TextValue $ws.Range("D11") "0.444"
$ws.Range("E11").Value = "  +1.01%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  -1.03%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "4.290.50"
$ws.Range("E13").Value = "  -0.59%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "32.45"
$ws.Range("E14").Value = "  +0.26%  "

# Row 15 - WrappedEther
Set-TextValue $ws.Range("D15") "3.680.66"
$ws.Range("E15").Value = "  -0.42%  "

# Row 16 - WrappedBTC
Set-TextValue $ws.Range("D16") "69.329.40"
$ws.Range("E16").Value = "  -0.15%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  -0.03%  "

# Row 18 - Chainlink
$ws.Range("E18").Value = "  -1.19%  "

# Row 19 - Polkadot
Set-TextValue $ws.Range("D19") "6.45"
$ws.Range("E19").Value = "  +0.00%  "

# Row 20 - BitcoinCash
Set-TextValue $ws.Range("D20") "465.82"
$ws.Range("E20").Value = "  -1.32%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "9.98"
$ws.Range("E21").Value = "  +0.83%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  -1.05%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "79.26"
$ws.Range("E23").Value = "  -1.09%  "

# Row 24 - WrappedeETH
Set-TextValue $ws.Range("D24") "3.818.31"
$ws.Range("E24").Value = "  -0.58%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.06%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  +0.01%  "

# Row 27 - InternetComputer(DFINITY) -- unchanged

# Row 28 - RenderToken
Set-TextValue $ws.Range("D28") "8.93"
$ws.Range("E28").Value = "  -2.19%  "

# Row 29 - PancakeSwap
Set-TextValue $ws.Range("D29") "2.62"
$ws.Range("E29").Value = "  -3.30%  "

# Row 30 - Fetch.AI
$ws.Range("E30").Value = "  -2.63%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("E31").Value = "  -0.19%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  -0.60%  "

# Row 33 - EthereumClassic
Set-TextValue $ws.Range("D33") "26.86"
$ws.Range("E33").Value = "  -0.31%  "

# Row 34 - NEARProtocol
$ws.Range("E34").Value = "  -2.65%  "

# Row 35 - RenzoRestakedETH
Set-TextValue $ws.Range("D35") "3.664.41"
$ws.Range("E35").Value = "  -0.45%  "

# Row 36 - Kaspa
Set-TextValue $ws.Range("D36") "0.161"
$ws.Range("E36").Value = "  -0.27%  "

# Row 37 - Aptos
Set-TextValue $ws.Range("D37") "8.35"
$ws.Range("E37").Value = "  -0.50%  "

# Row 38 - USDe -- unchanged

# Row 39 - Monero
Set-TextValue $ws.Range("D39") "178.77"
$ws.Range("E39").Value = "  +5.21%  "

# Row 40 - Filecoin
Set-TextValue $ws.Range("D40") "5.85"
$ws.Range("E40").Value = "  -6.90%  "

# Row 41 - FirstDigitalUSD
Set-TextValue $ws.Range("D41") "0.999"
$ws.Range("E41").Value = "  -0.08%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  -2.71%  "

# Row 43 - Hedera
Set-TextValue $ws.Range("D43") "0.0890"
$ws.Range("E43").Value = "  -1.86%  "

# Row 44 - Mantle
Set-TextValue $ws.Range("D44") "0.925"
$ws.Range("E44").Value = "  -1.82%  "

# Row 45 - OKB
Set-TextValue $ws.Range("D45") "46.74"
$ws.Range("E45").Value = "  -0.77%  "

# Row 46 - dogwifhat
$ws.Range("E46").Value = "  -0.82%  "

# Row 47 - InjectiveProtocol
Set-TextValue $ws.Range("D47") "27.86"
$ws.Range("E47").Value = "  -3.78%  "

# Row 48 - SuiNetwork
$ws.Range("E48").Value = "  -4.00%  "

# Row 49 - Cosmos
Set-TextValue $ws.Range("D49") "7.78"
$ws.Range("E49").Value = "  -0.86%  "

# Row 50 - FLOKI
$ws.Range("E50").Value = "  -5.31%  "
